$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 23-45 are brand new table rows (the original file had only a leftover
# empty, underlined placeholder cell at B25). Clear any inherited formatting
# over the whole target range first so the new rows don't pick up the old
# "end of table" underline style.
$ws.Range("A23:C45").ClearFormats() | Out-Null

$ws.Range("A23").Value = "UI.Prizes"
$ws.Range("B23").Value = "Prizes:"
$ws.Range("C23").Value = "Recompensas:"
$ws.Range("A24").Value = "Common.Money"
$ws.Range("B24").Value = "Gold"
$ws.Range("C24").Value = "Ouro"
$ws.Range("A25").Value = "Common.HealthPotion"
$ws.Range("B25").Value = "Health Potion"
$ws.Range("C25").Value = "Poção de Vida"
$ws.Range("A26").Value = "Common.MagicPotion"
$ws.Range("B26").Value = "Magic Potion"
$ws.Range("C26").Value = "Poção de Magia"
$ws.Range("A27").Value = "Common.Bomb"
$ws.Range("B27").Value = "Bomb"
$ws.Range("C27").Value = "Bomba"
$ws.Range("A28").Value = "Mission.Title01"
$ws.Range("B28").Value = "Kill Some Worms"
$ws.Range("C28").Value = "Mate algumas Minhocas"
$ws.Range("A29").Value = "Mission.Description01"
$ws.Range("B29").Value = "Worms Attacked the Forest. Kill them."
$ws.Range("C29").Value = "Minhocas atacaram a floresta. Mate-as."
$ws.Range("A30").Value = "Mission.Title02"
$ws.Range("B30").Value = "Get your bow!"
$ws.Range("C30").Value = "Obtenha seu arco!"
$ws.Range("A31").Value = "Mission.Description02"
$ws.Range("B31").Value = "Get your bow and kill a few long ranged enemies."
$ws.Range("C31").Value = "Obtenha seu arco e mate alguns inimigos com ataques de longo alcance."
$ws.Range("A32").Value = "Mission.Title03"
$ws.Range("B32").Value = "Bow Practice"
$ws.Range("C32").Value = "Prática de Arco e Flecha"
$ws.Range("A33").Value = "Mission.Description03"
$ws.Range("B33").Value = "Kill enemies and solve puzzles with the bow."
$ws.Range("C33").Value = "Mate inimigos e resolva quebra-cabeças com seu arco."
$ws.Range("A34").Value = "Mission.Title04"
$ws.Range("B34").Value = "Learn some Magic!"
$ws.Range("C34").Value = "Aprenda magia!"
$ws.Range("A35").Value = "Mission.Description04"
$ws.Range("B35").Value = "Get the fireball spell on the catacombs."
$ws.Range("C35").Value = "Obtenha a magia da bola de fogo nas catacombas."
$ws.Range("A36").Value = "Mission.Title05"
$ws.Range("B36").Value = "Do some magic!"
$ws.Range("C36").Value = "Pratique magia!"
$ws.Range("A37").Value = "Mission.Description05"
$ws.Range("B37").Value = "Explore a dark dungeon and kill some enemies with your magic."
$ws.Range("C37").Value = "Explore um calabouço escuro e mate alguns inimigos com sua magia."
$ws.Range("A38").Value = "Mission.Title06"
$ws.Range("B38").Value = "Haunted Boomerang."
$ws.Range("C38").Value = "Bumerangue mal-assombrado."
$ws.Range("A39").Value = "Mission.Description06"
$ws.Range("B39").Value = "People say the dark forest guards a treasure. Investigate it!"
$ws.Range("C39").Value = "Rumores dizem que a floresta negra guarda um grande tesouro. Investigue!"
$ws.Range("A40").Value = "Mission.Title07"
$ws.Range("B40").Value = "Boomerang practice"
$ws.Range("C40").Value = "Prática de Bumerangue"
$ws.Range("A41").Value = "Mission.Description07"
$ws.Range("B41").Value = "Kill waves after waves of enemies with the help of your boomerang!"
$ws.Range("C41").Value = "Mate hordas de inimigos com a ajuda do seu bumerangue!"
$ws.Range("A42").Value = "Mission.Title08"
$ws.Range("B42").Value = "The Royal Crypts"
$ws.Range("C42").Value = "Criptas Reais"
$ws.Range("A43").Value = "Mission.Description08"
$ws.Range("B43").Value = "Find the secret of the abandoned royal crypts."
$ws.Range("C43").Value = "Ache o segredo das criptas reais abandonadas."
$ws.Range("A44").Value = "Mission.Title09"
$ws.Range("B44").Value = "Get the crown!"
$ws.Range("C44").Value = "Recupere a coroa!"
$ws.Range("A45").Value = "Mission.Description09"
$ws.Range("B45").Value = "Kill the monsters and recover the crown!"
$ws.Range("C45").Value = "Mate os monstros e recupere a coroa!"

# Resize columns B and C to fit the new (longer) localized text, matching
# the widths recorded after the author auto-fit the columns.
$ws.Columns.Item(2).ColumnWidth = 62.285714285714285
$ws.Columns.Item(3).ColumnWidth = 65.42857142857143

# Restore the selection to the last edited cell, as in the saved workbook.
$ws.Range("C45").Select() | Out-Null
